# Apply "repull data" edits to column F (dSF) values.
# The diff only changes values in column F for a specific set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -1
    8  = -3
    16 = -6
    17 = 0
    19 = -3
    25 = 2
    26 = 5
    30 = -3
    31 = -8
    32 = -5
    39 = -7
    40 = 3
    41 = 1
    43 = -3
    46 = 1
    48 = 1
    51 = 9
    54 = -1
    56 = -1
    59 = -3
    62 = 3
    65 = -6
    69 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
